# daily auto push: 2025-10-12 09:24 UTC
# Append a new data row (row 94) to Sheet1, mirroring the existing rows'
# layout: A=date (text), B=weekday (text), C=hour (number), D=ranking (number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 94

# Force column A to stay a plain text value ("2025/10/12") instead of being
# auto-converted to a date serial number, then strip the formatting change
# so the cell keeps using the sheet's default (unstyled) cell format, just
# like the rest of the data rows.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/12"
$cellA.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "日"
$ws.Cells.Item($newRow, 3).Value = 18
$ws.Cells.Item($newRow, 4).Value = 201
